$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Nppc"
$ws.Cells.Item(2, 3).Value = "Npr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.647218666666667
$ws.Cells.Item(2, 8).Value = 7.941656
$ws.Cells.Item(2, 9).Value = 0.4640059894538357
$ws.Cells.Item(2, 10).Value = 0.4640059894538356
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 7.746472333333334
$ws.Cells.Item(2, 14).Value = 23.239417
$ws.Cells.Item(2, 15).Value = 0.2566490502786115
$ws.Cells.Item(2, 16).Value = 0.2566490502786115
$ws.Cells.Item(2, 17).Value = 20.50660616161689
$ws.Cells.Item(2, 18).Value = 184.559455454552
$ws.Cells.Item(2, 19).Value = 0.1190866965169144
$ws.Cells.Item(2, 20).Value = 0.1190866965169144
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Nppc"
$ws.Cells.Item(3, 3).Value = "Npr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.647218666666667
$ws.Cells.Item(3, 8).Value = 7.941656
$ws.Cells.Item(3, 9).Value = 0.4640059894538357
$ws.Cells.Item(3, 10).Value = 0.4640059894538356
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 19.17296733333333
$ws.Cells.Item(3, 14).Value = 57.518902
$ws.Cells.Item(3, 15).Value = 0.6352212523820423
$ws.Cells.Item(3, 16).Value = 0.6352212523820423
$ws.Cells.Item(3, 17).Value = 50.75503702019022
$ws.Cells.Item(3, 18).Value = 456.795333181712
$ws.Cells.Item(3, 19).Value = 0.2947464657336342
$ws.Cells.Item(3, 20).Value = 0.2947464657336342
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Nppc"
$ws.Cells.Item(4, 3).Value = "Npr2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.647218666666667
$ws.Cells.Item(4, 8).Value = 7.941656
$ws.Cells.Item(4, 9).Value = 0.4640059894538357
$ws.Cells.Item(4, 10).Value = 0.4640059894538356
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.263693
$ws.Cells.Item(4, 14).Value = 9.791079
$ws.Cells.Item(4, 15).Value = 0.1081296973393462
$ws.Cells.Item(4, 16).Value = 0.1081296973393462
$ws.Cells.Item(4, 17).Value = 8.639709031869334
$ws.Cells.Item(4, 18).Value = 77.757381286824
$ws.Cells.Item(4, 19).Value = 0.05017282720328712
$ws.Cells.Item(4, 20).Value = 0.05017282720328711
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Nppc"
$ws.Cells.Item(5, 3).Value = "Npr2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.057920333333333
$ws.Cells.Item(5, 8).Value = 9.173760999999999
$ws.Cells.Item(5, 9).Value = 0.5359940105461642
$ws.Cells.Item(5, 10).Value = 0.5359940105461642
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 7.746472333333334
$ws.Cells.Item(5, 14).Value = 23.239417
$ws.Cells.Item(5, 15).Value = 0.2566490502786115
$ws.Cells.Item(5, 16).Value = 0.2566490502786115
$ws.Cells.Item(5, 17).Value = 23.68809525970411
$ws.Cells.Item(5, 18).Value = 213.192857337337
$ws.Cells.Item(5, 19).Value = 0.1375623537616971
$ws.Cells.Item(5, 20).Value = 0.1375623537616971
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Nppc"
$ws.Cells.Item(6, 3).Value = "Npr2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 3.057920333333333
$ws.Cells.Item(6, 8).Value = 9.173760999999999
$ws.Cells.Item(6, 9).Value = 0.5359940105461642
$ws.Cells.Item(6, 10).Value = 0.5359940105461642
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 19.17296733333333
$ws.Cells.Item(6, 14).Value = 57.518902
$ws.Cells.Item(6, 15).Value = 0.6352212523820423
$ws.Cells.Item(6, 16).Value = 0.6352212523820423
$ws.Cells.Item(6, 17).Value = 58.62940665893577
$ws.Cells.Item(6, 18).Value = 527.6646599304219
$ws.Cells.Item(6, 19).Value = 0.340474786648408
$ws.Cells.Item(6, 20).Value = 0.340474786648408
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Nppc"
$ws.Cells.Item(7, 3).Value = "Npr2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 3.057920333333333
$ws.Cells.Item(7, 8).Value = 9.173760999999999
$ws.Cells.Item(7, 9).Value = 0.5359940105461642
$ws.Cells.Item(7, 10).Value = 0.5359940105461642
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.263693
$ws.Cells.Item(7, 14).Value = 9.791079
$ws.Cells.Item(7, 15).Value = 0.1081296973393462
$ws.Cells.Item(7, 16).Value = 0.1081296973393462
$ws.Cells.Item(7, 17).Value = 9.980113186457666
$ws.Cells.Item(7, 18).Value = 89.82101867811899
$ws.Cells.Item(7, 19).Value = 0.05795687013605907
$ws.Cells.Item(7, 20).Value = 0.05795687013605907
